$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h-volume-change (E) cells; three coin pairs
# (rows 40/41, 44/45, 50/51) also swap rank position (B/C/D/E all change).

$ws.Range("D2").Value = '58.291.20'
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").Value = '2.595.82'
$ws.Range("E3").Value = '  -0.81%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'521.76"
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").Value = "'143.87"
$ws.Range("E6").Value = '  +0.81%  '

$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").Value = '2.615.41'
$ws.Range("E9").Value = '  -0.39%  '

$ws.Range("D10").Value = "'6.51"
$ws.Range("E10").Value = '  -1.17%  '

$ws.Range("E11").Value = '  -1.55%  '

$ws.Range("E12").Value = '  +0.08%  '

$ws.Range("E13").Value = '  -0.36%  '

$ws.Range("D14").Value = '3.053.10'
$ws.Range("E14").Value = '  -0.81%  '

$ws.Range("D15").Value = '58.242.32'
$ws.Range("E15").Value = '  -0.13%  '

$ws.Range("D16").Value = "'20.41"
$ws.Range("E16").Value = '  -2.69%  '

$ws.Range("E17").Value = '  -1.64%  '

$ws.Range("D18").Value = '2.604.66'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").Value = "'339.58"
$ws.Range("E19").Value = '  +0.93%  '

$ws.Range("E20").Value = '  -0.99%  '

$ws.Range("D21").Value = "'10.26"
$ws.Range("E21").Value = '  -1.25%  '

$ws.Range("D22").Value = "'6.42"
$ws.Range("E22").Value = '  +2.13%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").Value = "'65.36"
$ws.Range("E24").Value = '  +1.52%  '

$ws.Range("D25").Value = "'0.168"
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("E26").Value = '  -2.80%  '

$ws.Range("D27").Value = '2.714.92'
$ws.Range("E27").Value = '  -0.65%  '

$ws.Range("D28").Value = "'0.996"
$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("D29").Value = "'7.02"
$ws.Range("E29").Value = '  -1.50%  '

$ws.Range("E30").Value = '  -5.16%  '

$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("E32").Value = '  -5.53%  '

$ws.Range("E33").Value = '  -0.43%  '

$ws.Range("D34").Value = "'18.79"
$ws.Range("E34").Value = '  +0.09%  '

$ws.Range("D35").Value = "'149.89"
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  -2.33%  '

$ws.Range("D37").Value = "'1.13"
$ws.Range("E37").Value = '  -4.34%  '

$ws.Range("D38").Value = "'0.873"
$ws.Range("E38").Value = '  -1.21%  '

$ws.Range("D39").Value = "'0.872"
$ws.Range("E39").Value = '  +1.84%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = "'1.46"
$ws.Range("E40").Value = '  +2.05%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").Value = "'35.97"
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("D42").Value = "'3.53"
$ws.Range("E42").Value = '  -2.77%  '

$ws.Range("E43").Value = '  -0.36%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = "'0.602"
$ws.Range("E44").Value = '  -0.23%  '

$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").Value = "'270.84"
$ws.Range("E45").Value = '  +0.62%  '

$ws.Range("D46").Value = "'0.0959"
$ws.Range("E46").Value = '  -0.91%  '

$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").Value = "'18.79"
$ws.Range("E48").Value = '  -1.83%  '

$ws.Range("D49").Value = "'0.0524"
$ws.Range("E49").Value = '  -1.66%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '1.972.00'
$ws.Range("E50").Value = '  -2.99%  '

$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'18.72"
$ws.Range("E51").Value = '  +2.40%  '
